# Scheduled market-data refresh: update currentAveragePrice* / Leve price & profit
# columns (H:N) across several sheets with newly-fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 5375.4346
$ws.Range("I80").Value = 312.30768
$ws.Range("J80").Value = 11957.5
$ws.Range("K80").Value = 936.92304
$ws.Range("L80").Value = 35872.5
$ws.Range("M80").Value = 61.07695999999999
$ws.Range("N80").Value = -37868.5
$ws.Range("H83").Value = 5375.4346
$ws.Range("I83").Value = 312.30768
$ws.Range("J83").Value = 11957.5
$ws.Range("K83").Value = 2810.76912
$ws.Range("L83").Value = 107617.5
$ws.Range("M83").Value = 2181.23088
$ws.Range("N83").Value = -117601.5
$ws.Range("H86").Value = 2617.25
$ws.Range("I86").Value = 1850.75
$ws.Range("J86").Value = 3000.5
$ws.Range("K86").Value = 1850.75
$ws.Range("L86").Value = 3000.5
$ws.Range("M86").Value = -727.75
$ws.Range("N86").Value = -5246.5
$ws.Range("H88").Value = 750.4
$ws.Range("I88").Value = 885.1429000000001
$ws.Range("J88").Value = 436
$ws.Range("K88").Value = 885.1429000000001
$ws.Range("L88").Value = 436
$ws.Range("M88").Value = -479.1429000000001
$ws.Range("N88").Value = -1248
$ws.Range("H89").Value = 2617.25
$ws.Range("I89").Value = 1850.75
$ws.Range("J89").Value = 3000.5
$ws.Range("K89").Value = 9253.75
$ws.Range("L89").Value = 15002.5
$ws.Range("M89").Value = -3637.75
$ws.Range("N89").Value = -26234.5
$ws.Range("H91").Value = 750.4
$ws.Range("I91").Value = 885.1429000000001
$ws.Range("J91").Value = 436
$ws.Range("K91").Value = 885.1429000000001
$ws.Range("L91").Value = 436
$ws.Range("M91").Value = 518.8570999999999
$ws.Range("N91").Value = -3244
$ws.Range("H112").Value = 2551.5715
$ws.Range("J112").Value = 2551.5715
$ws.Range("L112").Value = 7654.7145
$ws.Range("N112").Value = -9870.7145
$ws.Range("H127").Value = 647.625
$ws.Range("I127").Value = 407.0909
$ws.Range("J127").Value = 1176.8
$ws.Range("K127").Value = 1221.2727
$ws.Range("L127").Value = 3530.4
$ws.Range("M127").Value = 3738.7273
$ws.Range("N127").Value = -13450.4
$ws.Range("H129").Value = 874.86487
$ws.Range("I129").Value = 340
$ws.Range("J129").Value = 958.4375
$ws.Range("K129").Value = 1020
$ws.Range("L129").Value = 2875.3125
$ws.Range("M129").Value = 3980
$ws.Range("N129").Value = -12875.3125
$ws.Range("H137").Value = 3934.7917
$ws.Range("I137").Value = 2252.125
$ws.Range("J137").Value = 7300.125
$ws.Range("K137").Value = 6756.375
$ws.Range("L137").Value = 21900.375
$ws.Range("M137").Value = -4206.375
$ws.Range("N137").Value = -27000.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6466.3555
$ws.Range("I32").Value = 5310.405
$ws.Range("J32").Value = 22649.666
$ws.Range("K32").Value = 5310.405
$ws.Range("L32").Value = 22649.666
$ws.Range("M32").Value = -5023.405
$ws.Range("N32").Value = -23223.666
$ws.Range("H74").Value = 2598.762
$ws.Range("I74").Value = 2877.3333
$ws.Range("J74").Value = 1902.3334
$ws.Range("K74").Value = 2877.3333
$ws.Range("L74").Value = 1902.3334
$ws.Range("M74").Value = -2003.3333
$ws.Range("N74").Value = -3650.3334
$ws.Range("H77").Value = 2598.762
$ws.Range("I77").Value = 2877.3333
$ws.Range("J77").Value = 1902.3334
$ws.Range("K77").Value = 14386.6665
$ws.Range("L77").Value = 9511.666999999999
$ws.Range("M77").Value = -10018.6665
$ws.Range("N77").Value = -18247.667
$ws.Range("H88").Value = 1851.3
$ws.Range("I88").Value = 1751
$ws.Range("J88").Value = 2001.75
$ws.Range("K88").Value = 1751
$ws.Range("L88").Value = 2001.75
$ws.Range("M88").Value = -1345
$ws.Range("N88").Value = -2813.75
$ws.Range("H91").Value = 1851.3
$ws.Range("I91").Value = 1751
$ws.Range("J91").Value = 2001.75
$ws.Range("K91").Value = 1751
$ws.Range("L91").Value = 2001.75
$ws.Range("M91").Value = -347
$ws.Range("N91").Value = -4809.75
$ws.Range("H113").Value = 74099
$ws.Range("J113").Value = 74099
$ws.Range("L113").Value = 74099
$ws.Range("N113").Value = -82777
$ws.Range("H118").Value = 47000
$ws.Range("J118").Value = 47000
$ws.Range("L118").Value = 47000
$ws.Range("N118").Value = -50314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 24896.87
$ws.Range("J112").Value = 24896.87
$ws.Range("L112").Value = 24896.87
$ws.Range("N112").Value = -27850.87
$ws.Range("H134").Value = 2614.7222
$ws.Range("I134").Value = 1737.5834
$ws.Range("J134").Value = 4369
$ws.Range("K134").Value = 5212.7502
$ws.Range("L134").Value = 13107
$ws.Range("M134").Value = -2677.7502
$ws.Range("N134").Value = -18177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5697.7812
$ws.Range("I31").Value = 6879.0527
$ws.Range("J31").Value = 3971.3076
$ws.Range("K31").Value = 6879.0527
$ws.Range("L31").Value = 3971.3076
$ws.Range("M31").Value = -6584.0527
$ws.Range("N31").Value = -4561.3076
$ws.Range("H34").Value = 5697.7812
$ws.Range("I34").Value = 6879.0527
$ws.Range("J34").Value = 3971.3076
$ws.Range("K34").Value = 6879.0527
$ws.Range("L34").Value = 3971.3076
$ws.Range("M34").Value = -6677.0527
$ws.Range("N34").Value = -4375.3076
$ws.Range("H86").Value = 4505.35
$ws.Range("I86").Value = 4507.643
$ws.Range("K86").Value = 4507.643
$ws.Range("M86").Value = -3384.643
$ws.Range("H89").Value = 4505.35
$ws.Range("I89").Value = 4507.643
$ws.Range("K89").Value = 22538.215
$ws.Range("M89").Value = -16922.215

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2792.3
$ws.Range("I132").Value = 2145.2307
$ws.Range("J132").Value = 3994
$ws.Range("K132").Value = 6435.6921
$ws.Range("L132").Value = 11982
$ws.Range("M132").Value = -3905.6921
$ws.Range("N132").Value = -17042

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1850
$ws.Range("J68").Value = 2533.3333
$ws.Range("L68").Value = 2533.3333
$ws.Range("N68").Value = -4031.3333
$ws.Range("H71").Value = 1850
$ws.Range("J71").Value = 2533.3333
$ws.Range("L71").Value = 12666.6665
$ws.Range("N71").Value = -20154.6665
$ws.Range("H132").Value = 3291.413
$ws.Range("I132").Value = 2994.0667
$ws.Range("J132").Value = 3848.9375
$ws.Range("K132").Value = 8982.2001
$ws.Range("L132").Value = 11546.8125
$ws.Range("M132").Value = -6452.2001
$ws.Range("N132").Value = -16606.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3799.2
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 3999.111
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 7998.222
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -10120.222
$ws.Range("H84").Value = 3799.2
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 3999.111
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 39991.11
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -50599.11
$ws.Range("H111").Value = 49643.5
$ws.Range("J111").Value = 49643.5
$ws.Range("L111").Value = 49643.5
$ws.Range("N111").Value = -57823.5
